$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style (date/number formats) from the last existing row (621)
# down across the new block of rows (622:651) in one broadcast paste.
$ws.Range("A621:C621").Copy()
$ws.Range("A622:C651").PasteSpecial(-4122)

# Populate the new daily UF / IVP values (10-Sep-2020 .. 09-Oct-2020).
$ws.Range("A622").Value = 44084
$ws.Range("B622").Value = 28688.73
$ws.Range("C622").Value = 30055.57
$ws.Range("A623").Value = 44085
$ws.Range("B623").Value = 28689.68
$ws.Range("C623").Value = 30056.23
$ws.Range("A624").Value = 44086
$ws.Range("B624").Value = 28690.639999999999
$ws.Range("C624").Value = 30056.880000000001
$ws.Range("A625").Value = 44087
$ws.Range("B625").Value = 28691.59
$ws.Range("C625").Value = 30057.54
$ws.Range("A626").Value = 44088
$ws.Range("B626").Value = 28692.55
$ws.Range("C626").Value = 30058.19
$ws.Range("A627").Value = 44089
$ws.Range("B627").Value = 28693.51
$ws.Range("C627").Value = 30058.84
$ws.Range("A628").Value = 44090
$ws.Range("B628").Value = 28694.46
$ws.Range("C628").Value = 30059.5
$ws.Range("A629").Value = 44091
$ws.Range("B629").Value = 28695.42
$ws.Range("C629").Value = 30060.15
$ws.Range("A630").Value = 44092
$ws.Range("B630").Value = 28696.37
$ws.Range("C630").Value = 30060.81
$ws.Range("A631").Value = 44093
$ws.Range("B631").Value = 28697.33
$ws.Range("C631").Value = 30061.46
$ws.Range("A632").Value = 44094
$ws.Range("B632").Value = 28698.29
$ws.Range("C632").Value = 30062.11
$ws.Range("A633").Value = 44095
$ws.Range("B633").Value = 28699.24
$ws.Range("C633").Value = 30062.77
$ws.Range("A634").Value = 44096
$ws.Range("B634").Value = 28700.2
$ws.Range("C634").Value = 30063.42
$ws.Range("A635").Value = 44097
$ws.Range("B635").Value = 28701.15
$ws.Range("C635").Value = 30064.080000000002
$ws.Range("A636").Value = 44098
$ws.Range("B636").Value = 28702.11
$ws.Range("C636").Value = 30064.73
$ws.Range("A637").Value = 44099
$ws.Range("B637").Value = 28703.07
$ws.Range("C637").Value = 30065.38
$ws.Range("A638").Value = 44100
$ws.Range("B638").Value = 28704.02
$ws.Range("C638").Value = 30066.04
$ws.Range("A639").Value = 44101
$ws.Range("B639").Value = 28704.98
$ws.Range("C639").Value = 30066.69
$ws.Range("A640").Value = 44102
$ws.Range("B640").Value = 28705.94
$ws.Range("C640").Value = 30067.35
$ws.Range("A641").Value = 44103
$ws.Range("B641").Value = 28706.89
$ws.Range("C641").Value = 30068
$ws.Range("A642").Value = 44104
$ws.Range("B642").Value = 28707.85
$ws.Range("C642").Value = 30068.65
$ws.Range("A643").Value = 44105
$ws.Range("B643").Value = 28708.799999999999
$ws.Range("C643").Value = 30069.31
$ws.Range("A644").Value = 44106
$ws.Range("B644").Value = 28709.759999999998
$ws.Range("C644").Value = 30069.96
$ws.Range("A645").Value = 44107
$ws.Range("B645").Value = 28710.720000000001
$ws.Range("C645").Value = 30070.62
$ws.Range("A646").Value = 44108
$ws.Range("B646").Value = 28711.67
$ws.Range("C646").Value = 30071.27
$ws.Range("A647").Value = 44109
$ws.Range("B647").Value = 28712.63
$ws.Range("C647").Value = 30071.93
$ws.Range("A648").Value = 44110
$ws.Range("B648").Value = 28713.59
$ws.Range("C648").Value = 30072.58
$ws.Range("A649").Value = 44111
$ws.Range("B649").Value = 28714.54
$ws.Range("C649").Value = 30073.23
$ws.Range("A650").Value = 44112
$ws.Range("B650").Value = 28715.5
$ws.Range("C650").Value = 30073.89
$ws.Range("A651").Value = 44113
$ws.Range("B651").Value = 28716.46
$ws.Range("C651").Value = 30074.54

# Extend the workbook-level defined name to cover the new data range.
$nm = $wb.Names.Item(1)
$nm.RefersTo = "=UF_IVP_DIARIO!`$A`$1:`$C`$651"

# Update the view state (scroll position / active selection) to mirror
# where the author left the cursor after the update.
$win = $excel.ActiveWindow
$win.ScrollRow = 639
$win.ScrollColumn = 1
[void]$ws.Range("C669").Select()

Write-Output "UF_IVP_DIARIO updated through row 651"
